$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.877.58'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.892.23'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7752'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3142'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07407'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08148'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7657'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.480'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.89%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.863.61'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.52'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.222'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.98%  '
$ws.Range('D17').Value = '29.849.96'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.40'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007873'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9993'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.128'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.83%  '
$ws.Range('D23').Value = '2.120.45'
$ws.Range('E23').Value = '  -0.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1570'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.429'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.80'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.041'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.455'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.545'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.503'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05605'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.096'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7586'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9978'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.648'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01937'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.791'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('D41').Value = '1.148.00'
$ws.Range('E41').Value = '  +12.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.38'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4462'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.977'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8555'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.000'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.136'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '101.85'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.897'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.518'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.31%  '
